{"js": "// Append \": Y\" to the label text of specific table-cell variable rows.\n// Each target string is the exact, whole text of a <w:t> run that sits\n// alone in its paragraph (first column of the results table), so a plain\n// literal search + whole-range replace is safe and keeps run/paragraph\n// formatting untouched.\nconst targets = [\n  \"postop_MCS_dependence\",\n  \"postop_VA_ECMO\",\n  \"postop_BiVAD_dependence\",\n  \"postop_MCS_Impella5.5_DEPENDENT\",\n  \"postop_CRRT\",\n  \"postop_stroke\",\n  \"ACR_2R_or_greater\",\n  \"survival_90\",\n];\n\nconst body = context.document.body;\n\nfor (const target of targets) {\n  const results = body.search(target, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  for (const r of results.items) {\n    if (r.text === target) {\n      r.insertText(target + \": Y\", Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Append \": Y\" to the label text of specific table-cell variable rows.\n# Each target string is the exact, whole text of a single-run paragraph\n# (first column of the results table), so Find/Replace over the whole\n# document content is safe and keeps run/paragraph formatting untouched.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"postop_MCS_dependence\",\n    \"postop_VA_ECMO\",\n    \"postop_BiVAD_dependence\",\n    \"postop_MCS_Impella5.5_DEPENDENT\",\n    \"postop_CRRT\",\n    \"postop_stroke\",\n    \"ACR_2R_or_greater\",\n    \"survival_90\"\n)\n\nforeach ($target in $targets) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $target\n    $find.Replacement.Text = $target + \": Y\"\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($target, $true, $false, $false, $false, $false, $true, 0, $false, ($target + \": Y\"), 2) | Out-Null\n}\n"}
